$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (9) to hold the new
# "Distribution channel code" field; this shifts the existing
# "Logistic percentage" column from I to J.
$ws.Columns.Item(9).Insert()

# Header + values for the newly inserted column.
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Match the manually-resized column width Excel recorded for the
# new column (slightly narrower than the neighboring J column).
$ws.Columns.Item(9).ColumnWidth = 21.67

Write-Host "done"
